$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column C: replace the long math text with the PDF reference and
# drop the wrap-text formatting / huge row height that went with it.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(2, 3).PasteSpecial(-4122) | Out-Null  # xlPasteFormats -> default style
$ws.Cells.Item(2, 3).Value = "pdf/GTC1.pdf"
$ws.Rows.Item(2).AutoFit()

# Row 3: give the subject a LaTeX bold wrapper, and point the doc column at
# the same pdf file.
$ws.Cells.Item(3, 2).Value = "\(\textbf{Giới hạn của hàm số}\)"
$ws.Cells.Item(3, 3).Value = "pdf/GTC1.pdf"

# New row 4: another topic referencing the same pdf.
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "\(\textbf{Giới hạn lượng giác}\)"
$ws.Cells.Item(4, 3).Value = "pdf/GTC1.pdf"

$ws.Range("B4").Select() | Out-Null
